$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new BOM row (row 16): Common cathode Schottky diode (BAT160C) ---
# Values are entered in the order purpose -> link -> part# so the new
# shared-string table entries land in the same order as the target file:
#   59 = "Common cathode Shotkey", 60 = datasheet/part-link URL, 61 = "BAT160C"
$ws.Range("A16").Value = "Common cathode Shotkey"
$ws.Range("C16").Value = "https://www.digikey.ca/en/products/detail/nexperia-usa-inc/BAT160C-115/1232113"
$ws.Range("B16").Value = "BAT160C"

# The part-link cell for this row uses a bold, wrapped style (new font/style).
$ws.Range("C16").Font.Bold = $true
$ws.Range("C16").WrapText = $true

# Row height matches the other 75pt rows used for similar entries.
$ws.Rows.Item(16).RowHeight = 75

# --- Turn the existing LDO 5v->3v3 datasheet text (C2) into a real hyperlink ---
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.digikey.ca/en/products/detail/diodes-incorporated/AZ1117IH-3-3TRG1/5699672") | Out-Null
$ws.Range("C2").WrapText = $true

# --- Update the visible selection to match the new active cell ---
$ws.Activate()
$ws.Range("B16").Select() | Out-Null
